# cx_Freeze compile with database: death to Python2.7
#
# Applies the REV-log updates: rewrites the "REV 0" note, appends
# "REV 8"/"REV 9"/"REV x.x" rows to the revision sheet, and bumps the
# "Created" timestamp on the Parts - Consoles sheet.

$wb = $excel.ActiveWorkbook

# ---- Sheet "revision" -----------------------------------------------
$rev = $wb.Worksheets.Item("revision")

# Row 1: reword the REV 0 note, keep the same author id
$rev.Range("C1").Value = "Order from chaos"
$rev.Range("D1").Value = "id-0"

# Row 9 ("REV") becomes "REV 8"
$rev.Range("A9").Value = "REV 8"

# Row 10 used to hold the closing note with a live timestamp in column B
# and the old closing quote in column C; it becomes a normal text row
# for REV 9, and the timestamp + closing quote move down to a new row 11.
$rev.Range("A10").Value = "REV 9"
# the date goes in as plain text here (not a real date value), so force
# text and strip the date-time formatting the cell inherited
$rev.Range("B10").Value = "'20180307"
$rev.Range("B10").Style = "Normal"
$rev.Range("C10").Value = "Test build with cx_Freeze5.1.1"
$rev.Range("D10").Value = "id-0"

$rev.Range("A11").Value = "REV x.x"
$rev.Range("B11").Value = 43166.75549595193
$rev.Range("B11").NumberFormat = "yyyy-mm-dd h:mm:ss"
$rev.Range("C11").Value = "Fiat justitia ruat caelum "
$rev.Range("D11").Value = "id-0"

# Column widths were nudged when the table grew a row
$rev.Columns.Item(1).ColumnWidth = 10
$rev.Columns.Item(3).ColumnWidth = 46

# ---- Sheet "Parts - Consoles" ----------------------------------------
$parts = $wb.Worksheets.Item("Parts - Consoles")
$parts.Range("BA2").Value = 43166.75561099128
